# Fix the "Elastic network interface" textbox labels so they match the
# sizing/font used by the rest of the architecture diagram's callout labels
# (ext cy 430887 -> 461665 EMU, i.e. Height 33.92811pt -> 36.35158pt, and
# run font size 11pt -> 12pt).
#
# Slide 1: only the box anchored at off x="4612866" y="5125416" is updated.
# Slide 2: only the box anchored at off x="7884631" y="5125416" is updated.
# (The other "Elastic network interface" box on each slide is left as-is.)

$p = $ppt.ActivePresentation

$targetHeightPt = 36.35158   # rounds to exactly 461665 EMU
$targetFontSize = 12         # 1200 (centi-points)
$targetText = "Elastic network interface"

# (slide index, expected Left in points) for the single box to fix per slide
$fixups = @(
    @{ Slide = 1; Left = 363.2178 },
    @{ Slide = 2; Left = 620.8371 }
)

foreach ($fix in $fixups) {
    $s = $p.Slides.Item($fix.Slide)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $targetText) {
                if ([Math]::Abs($sh.Left - $fix.Left) -lt 0.01) {
                    $sh.Height = $targetHeightPt
                    $sh.TextFrame.TextRange.Font.Size = $targetFontSize
                }
            }
        }
    }
}
